$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_R_acc_LT"

$ws.Range("A2").Value = 80.301645338208402
$ws.Range("A3").Value = 82.17550274223035
$ws.Range("A4").Value = 83.135283363802557
$ws.Range("A5").Value = 87.477148080438766
$ws.Range("A6").Value = 87.659963436928706
$ws.Range("A7").Value = 87.340036563071294
$ws.Range("A8").Value = 78.244972577696529
$ws.Range("A9").Value = 78.427787934186483
$ws.Range("A10").Value = 77.970749542961599
$ws.Range("A11").Value = 73.49177330895796
$ws.Range("A12").Value = 74.177330895795251
$ws.Range("A13").Value = 79.524680073126149
$ws.Range("A14").Value = 77.559414990859239
$ws.Range("A15").Value = 77.925045703839118
$ws.Range("A16").Value = 77.513711151736743
$ws.Range("A17").Value = 76.873857404021933
$ws.Range("A18").Value = 79.707495429616088
$ws.Range("A19").Value = 84.643510054844612
$ws.Range("A20").Value = 87.659963436928706
$ws.Range("A21").Value = 85.557586837294338
$ws.Range("A22").Value = 87.477148080438766
$ws.Range("A23").Value = 79.204753199268737
$ws.Range("A24").Value = 81.946983546617915
$ws.Range("A25").Value = 81.261425959780624
$ws.Range("A26").Value = 78.793418647166362
$ws.Range("A27").Value = 78.199268738574034
$ws.Range("A28").Value = 79.296160877513714
$ws.Range("A29").Value = 80.712979890310791
$ws.Range("A30").Value = 79.57038391224863
$ws.Range("A31").Value = 79.204753199268737
$ws.Range("A32").Value = 88.482632541133455
$ws.Range("A33").Value = 89.899451553930533
$ws.Range("A34").Value = 91.133455210237656
$ws.Range("A35").Value = 81.809872029250457
$ws.Range("A36").Value = 79.75319926873857
$ws.Range("A37").Value = 66.453382084095068
$ws.Range("A38").Value = 84.3235831809872
$ws.Range("A39").Value = 76.279707495429619
$ws.Range("A40").Value = 79.616087751371111
$ws.Range("A41").Value = 78.610603290676423
$ws.Range("A42").Value = 78.427787934186483
$ws.Range("A43").Value = 78.976234003656316
$ws.Range("A44").Value = 78.244972577696529
$ws.Range("A45").Value = 78.976234003656316
$ws.Range("A46").Value = 78.884826325411339
$ws.Range("A47").Value = 78.564899451553927
$ws.Range("A48").Value = 74.771480804387565
$ws.Range("A49").Value = 79.981718464351005
